# Revert "Drop in files from RMI script"
#
# Re-introduces the "Texas Notes" worksheet (between "About" and "DR"),
# restores the discount-rate value on the DR sheet to 5.87% (VCE WISdom
# number) instead of the flat 3% RMI used, and sets the iterative
# calculation options / selections that went along with that edit.

$wb = $excel.ActiveWorkbook

# --- Enable iterative calculation (workbook previously had none) ---------
try {
    $excel.Iteration    = $true
    $excel.MaxIterations = 100
    $excel.MaxChange     = 0.00001
} catch {}

# --- Insert the "Texas Notes" sheet right before "DR" ---------------------
$drSheetForInsert = $wb.Worksheets.Item("DR")
$notesSheet = $wb.Worksheets.Add($drSheetForInsert)
$notesSheet.Name = "Texas Notes"

$notesSheet.Range("A1").Value = "updated to the VCE WISdom number"
$notesSheet.Range("A2").Value = 0.0587
$notesSheet.Range("A4").Value = "their feedback was the 3% was a bit low"

# --- Update the DR sheet's discount rate to match --------------------------
# (re-fetch by name: inserting the new sheet shifted DR's index)
$drSheet = $wb.Worksheets.Item("DR")
$drSheet.Range("B2").Value = 0.0587

# --- Restore/match selections on each sheet --------------------------------
$aboutSheet = $wb.Worksheets.Item("About")
[void]$aboutSheet.Activate()
[void]$aboutSheet.Range("C23").Select()

[void]$wb.Worksheets.Item("Texas Notes").Activate()
[void]$wb.Worksheets.Item("Texas Notes").Range("A5").Select()

[void]$wb.Worksheets.Item("DR").Activate()
[void]$wb.Worksheets.Item("DR").Range("B2").Select()

# "About" stays the active/visible tab
[void]$aboutSheet.Activate()
